$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated crypto data cell from the refreshed GitHub Actions run.
# Column D cells are forced to Text format first so that numeric-looking
# strings (e.g. "1.00", "12.60", "65.714.53") keep their exact display text
# instead of being auto-coerced into numbers by COM.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.714.53"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.653.43"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.88"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.80"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000198"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.71"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.130.31"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.563.19"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.632.45"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.46"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.20"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.02"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.66"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "533.39"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.44"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.37"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "155.75"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.37"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0607"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0255"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("E50").Value = "  +8.78%  "
$ws.Range("E51").Value = "  -2.50%  "
